$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (values must remain text, matching the
# original inline-string cell contents such as prices and percentages).
$updates = @{
    'D2' = '308.29'
    'E2' = '1.13%'
    'D3' = '36.38'
    'E3' = '1.41%'
    'E4' = '0.63%'
    'D5' = '0.08117'
    'E5' = '0.38%'
    'D6' = '2.055'
    'E6' = '7.56%'
    'D7' = '4.159'
    'E7' = '0.27%'
    'D8' = '7.859'
    'E8' = '0.19%'
    'E9' = '-0.46%'
    'D10' = '0.1424'
    'E10' = '12.97%'
    'D11' = '0.1927'
    'E11' = '0.72%'
    'D12' = '0.09098'
    'E12' = '-0.93%'
    'D13' = '0.03447'
    'E13' = '-1.38%'
    'D14' = '0.09910'
    'E14' = '-0.14%'
    'D15' = '0.001407'
    'E15' = '-0.44%'
    'D16' = '0.006271'
    'E16' = '-5.52%'
    'D17' = '3.842'
    'E17' = '6.15%'
    'D18' = '3.364'
    'E18' = '3.67%'
    'E19' = '-0.11%'
    'E20' = '-3.47%'
    'D21' = '4.800'
    'E21' = '-7.86%'
    'E22' = '-7.55%'
    'D23' = '0.04383'
    'E23' = '-0.56%'
    'D24' = '0.001234'
    'E24' = '-0.05%'
    'E25' = '4.10%'
    'E27' = '-0.10%'
    'D39' = '0.02028'
    'E39' = '3.56%'
    'D40' = '0.05159'
    'E40' = '-0.16%'
    'D41' = '0.007488'
    'E41' = '-1.41%'
    'D42' = '0.01010'
    'E42' = '0.46%'
    'D43' = '0.1372'
    'E43' = '0.33%'
    'D44' = '0.002130'
    'E44' = '1.32%'
    'D45' = '0.009965'
    'E45' = '-6.79%'
    'D46' = '0.00006313'
    'E46' = '-0.69%'
    'E47' = '-0.09%'
    'D48' = '64.86'
    'E48' = '-0.16%'
    'E49' = '-22.05%'
    'E50' = '-0.09%'
    'E51' = '-0.09%'
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    # Force text format so numeric-looking strings (e.g. "308.29" or
    # "1.13%") are stored as text instead of being parsed into numbers.
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
    # Reset the style back to the default so no stray style index is
    # introduced (the original cells used the default/general style).
    $range.Style = "Normal"
}
